$wb = $excel.ActiveWorkbook

# Overview sheet: mark the c7a2c94e row as handed back (matches zh-cn/de-de sheets)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: update status + handback datetime for c7a2c94e row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Handed back: in sync with en-US"
$wsZh.Range("G2").Value = "2016-02-18 10:33:02"
$wsZh.Range("G3").Value = "2016-02-18 10:33:02"

# de-de sheet: update status + handback datetime for c7a2c94e row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDe.Range("G2").Value = "2016-02-18 10:33:25"
$wsDe.Range("G3").Value = "2016-02-18 10:33:25"
